$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.052.96'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.327.60'
$ws.Range("E3").Value = '  +3.99%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '98.44'
$ws.Range("E5").Value = '  +4.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '271.16'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -2.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.77'
$ws.Range("E10").Value = '  -1.15%  '
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("E12").Value = '  -4.50%  '
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.668.11'
$ws.Range("E14").Value = '  +3.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.59'
$ws.Range("E15").Value = '  +1.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.883'
$ws.Range("E16").Value = '  +7.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.331.98'
$ws.Range("E17").Value = '  +4.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '44.006.92'
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("E19").Value = '  +4.64%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.41'
$ws.Range("E20").Value = '  +3.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.71'
$ws.Range("E21").Value = '  +3.99%  '
$ws.Range("E22").Value = '  -0.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '240.39'
$ws.Range("E23").Value = '  +3.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.31'
$ws.Range("E24").Value = '  +1.61%  '
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("E26").Value = '  +0.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.46'
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  -2.01%  '
$ws.Range("E29").Value = '  +2.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.58'
$ws.Range("E30").Value = '  -4.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.52'
$ws.Range("E31").Value = '  +7.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '176.26'
$ws.Range("E32").Value = '  +2.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0913'
$ws.Range("E33").Value = '  -0.38%  '
$ws.Range("E34").Value = '  +0.58%  '
$ws.Range("E35").Value = '  +2.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0364'
$ws.Range("E36").Value = '  +3.28%  '
$ws.Range("E37").Value = '  -2.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.48'
$ws.Range("E38").Value = '  +4.11%  '
$ws.Range("E39").Value = '  -4.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.45'
$ws.Range("E40").Value = '  +13.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.244'
$ws.Range("E41").Value = '  +8.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.44'
$ws.Range("E42").Value = '  +23.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.45'
$ws.Range("E43").Value = '  -3.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.89'
$ws.Range("E44").Value = '  -0.69%  '
$ws.Range("E45").Value = '  +8.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.35'
$ws.Range("E46").Value = '  -1.20%  '
$ws.Range("E47").Value = '  +4.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '100.53'
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("E49").Value = '  +0.87%  '
$ws.Range("E50").Value = '  +15.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.554.93'
$ws.Range("E51").Value = '  +4.03%  '
